$d = $word.ActiveDocument
$sec = $d.Sections.First

# --- Footer1: rename inline picture "image2.png" -> "image1.png" (docPr id=1, descr=Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png) ---
$hf = $sec.Footers.Item(1)
if ($hf.Exists) {
    $shp = $null
    for ($i = 1; $i -le $hf.Range.InlineShapes.Count; $i++) {
        $cand = $hf.Range.InlineShapes.Item($i)
        if (($cand.AlternativeText -eq 'Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png') -and ($cand.Width -ne $null)) {
            $shp = $cand
        }
    }
    if ($shp -ne $null) {
        $rng = $shp.Range
        $shp.Delete()
        [void]$rng.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:b w:val="1"/></w:rPr><w:drawing><wp:inline distB="0" distT="0" distL="0" distR="0"><wp:extent cx="952500" cy="285750"/><wp:effectExtent b="0" l="0" r="0" t="0"/><wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="1" name="image1.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"/><pic:cNvPicPr preferRelativeResize="0"/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId1"/><a:srcRect b="0" l="0" r="0" t="0"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="952500" cy="285750"/></a:xfrm><a:prstGeom prst="rect"/><a:ln/></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>')
    }
}

# --- Footer2: rename inline picture "image2.png" -> "image1.png" (docPr id=2, descr=Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png) ---
$hf = $sec.Footers.Item(2)
if ($hf.Exists) {
    $shp = $null
    for ($i = 1; $i -le $hf.Range.InlineShapes.Count; $i++) {
        $cand = $hf.Range.InlineShapes.Item($i)
        if (($cand.AlternativeText -eq 'Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png') -and ($cand.Width -ne $null)) {
            $shp = $cand
        }
    }
    if ($shp -ne $null) {
        $rng = $shp.Range
        $shp.Delete()
        [void]$rng.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:rFonts w:ascii="Verdana" w:cs="Verdana" w:eastAsia="Verdana" w:hAnsi="Verdana"/><w:b w:val="1"/><w:i w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="none"/><w:shd w:fill="auto" w:val="clear"/><w:vertAlign w:val="baseline"/></w:rPr><w:drawing><wp:inline distB="0" distT="0" distL="0" distR="0"><wp:extent cx="952500" cy="285750"/><wp:effectExtent b="0" l="0" r="0" t="0"/><wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image1.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"/><pic:cNvPicPr preferRelativeResize="0"/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId1"/><a:srcRect b="0" l="0" r="0" t="0"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="952500" cy="285750"/></a:xfrm><a:prstGeom prst="rect"/><a:ln/></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>')
    }
}

# --- Header2: rename inline picture "image1.jpg" -> "image2.jpg" (docPr id=3, descr=BTec_Logo-Orange) ---
$hf = $sec.Headers.Item(2)
if ($hf.Exists) {
    $shp = $null
    for ($i = 1; $i -le $hf.Range.InlineShapes.Count; $i++) {
        $cand = $hf.Range.InlineShapes.Item($i)
        if (($cand.AlternativeText -eq 'BTec_Logo-Orange') -and ($cand.Width -ne $null)) {
            $shp = $cand
        }
    }
    if ($shp -ne $null) {
        $rng = $shp.Range
        $shp.Delete()
        [void]$rng.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:rFonts w:ascii="Verdana" w:cs="Verdana" w:eastAsia="Verdana" w:hAnsi="Verdana"/><w:b w:val="0"/><w:i w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="none"/><w:shd w:fill="auto" w:val="clear"/><w:vertAlign w:val="baseline"/></w:rPr><w:drawing><wp:inline distB="0" distT="0" distL="0" distR="0"><wp:extent cx="914400" cy="277792"/><wp:effectExtent b="0" l="0" r="0" t="0"/><wp:docPr descr="BTec_Logo-Orange" id="3" name="image2.jpg"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image2.jpg"/><pic:cNvPicPr preferRelativeResize="0"/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId1"/><a:srcRect b="0" l="0" r="0" t="0"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="914400" cy="277792"/></a:xfrm><a:prstGeom prst="rect"/><a:ln/></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>')
    }
}
